$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing row (old row 24 / Requisitos detail moves up to row 23)
$ws.Rows.Item(24).Delete()

# Row 1
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

# Row 2
$ws.Range("B2").Value = "LOQ4081"
$ws.Range("C2").Value = "LOQ4081"

# Row 3
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Química Orgânica"
$ws.Range("C3").Value = " Química Orgânica"

# Row 4
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Organic Chemistry"
$ws.Range("C4").Value = "Organic Chemistry"

# Row 5
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

# Row 6
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

# Row 7
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

# Row 8
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Row 9
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EA-3"
$ws.Range("C9").Value = "EA-3"

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Range("C10").Value = "2346890 - Eliane Corrêa Pedrozo"

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Study of the classes of organic compounds most relevant from the point of view of its presence in the environment. Understand the relationship of molecular structures to physical properties and their chemical reactivity. Propose appropriate actions that can eliminate or mitigate harmful effects related to them in order to enable sustainable growth and development for future generations."
$ws.Range("C11").Value = "Study of the classes of organic compounds most relevant from the point of view of its presence in the environment. Understand the relationship of molecular structures to physical properties and their chemical reactivity. Propose appropriate actions that can eliminate or mitigate harmful effects related to them in order to enable sustainable growth and development for future generations."

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds."
$ws.Range("C14").Value = "Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds."
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers."
$ws.Range("C16").Value = "Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers."
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Range("C18").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas sobre os tópicos descritos no programa. Seminários e/ou projeto dirigido. Provas escritas."
$ws.Range("C19").Value = "Aulas expositivas sobre os tópicos descritos no programa. Seminários e/ou projeto dirigido. Provas escritas."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Avaliação individual escrita e apresentação de seminário e/ou projeto. Avaliações individuais correspondem a 70% da NF (nota final) e em grupo 30% da NF. Se NF for => 5,0: Aluno aprovado, condicionado à frequência mínima de 70%."
$ws.Range("C20").Value = "Avaliação individual escrita e apresentação de seminário e/ou projeto. Avaliações individuais correspondem a 70% da NF (nota final) e em grupo 30% da NF. Se NF for => 5,0: Aluno aprovado, condicionado à frequência mínima de 70%."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Prova escrita envolvendo todo o conteúdo da disciplina. MF = (NF+R) / 2 & => 5,0 Aprovado"
$ws.Range("C21").Value = "Prova escrita envolvendo todo o conteúdo da disciplina. MF = (NF+R) / 2 & => 5,0 Aprovado"
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# Row 23
$ws.Range("B23").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
"
$ws.Range("C23").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)
"
$ws.Range("A23").ClearContents()
$ws.Rows.Item(23).RowHeight = 30
